$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "8'. AHE cl interp buf uchar (32 bin)" row.
# This shifts the old row 11 ("9. AHE cl interp buf uchar bank" with its data)
# up into row 10, matching the target state exactly.
$ws.Rows.Item(10).Delete()

# Append a brand new row of results: "10. AHE cl interp buf uchar bank mod"
$ws.Range("A11").Value = "10. AHE cl interp buf uchar bank mod"
$ws.Range("B11").Value = 241.74392
$ws.Range("C11").Value = 399.85935999999998
$ws.Range("D11").Value = 399.85935999999998
$ws.Range("E11").Value = 1.94112

# C11/E11 pick up the "Times / explicit black" number format already used
# elsewhere in the sheet (e.g. D2), so copy that formatting over.
$ws.Range("D2").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("E11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Widen column A a bit to fit the longer label text.
$ws.Columns.Item(1).ColumnWidth = 34.83

# Leave the cursor/selection where the author last left it.
$ws.Range("C16").Select()
